$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 141, shifting existing rows 141-218 down to 144-221
$ws.Rows("141:143").Insert()

# Fill in the new rows 141-143 with data for date 44438 (2021-09-06)
$newRows = @(
    @{Row=141; D=44438; I="Especial"; J=500; K=16000; L=17000; M=16500; P=917},
    @{Row=142; D=44438; I="Primera";  J=400; K=14000; L=15000; M=14500; P=806},
    @{Row=143; D=44438; I="Segunda";  J=300; K=12000; L=12500; M=12250; P=681}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = 'Terminal La Palmera de La Serena'
    $ws.Cells.Item($row, 3).Value = 'Coquimbo'
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = 'Pepino dulce'
    $ws.Cells.Item($row, 8).Value = 'Cultivar IV Región'
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = '$/bandeja 18 kilos'
    $ws.Cells.Item($row, 15).Value = 'Provincia de Limarí'
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = 'Hortaliza'
}
